# Added normal price support. Need to add photos
#
# The "price" column (D) on the "Assemblies" sheet used to hold formatted
# text like "269 790 ₽". Convert every row's price into a plain numeric
# value (the ruble sign / thousands separators are dropped) so downstream
# consumers can treat it as a normal number.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assemblies")

# row -> numeric price (was a shared text string like "269 790 ₽")
$ws.Range("D2").Value = 269790
$ws.Range("D3").Value = 67900
$ws.Range("D4").Value = 25990
$ws.Range("D5").Value = 199500
$ws.Range("D6").Value = 116320
$ws.Range("D7").Value = 75960
$ws.Range("D8").Value = 99060

# Selection moved from G8 to D9 (next empty row in the price column).
$ws.Range("D9").Select()

# Window position bookkeeping (best effort; cosmetic only).
$excel.ActiveWindow.Left = 5910
$excel.ActiveWindow.Top = 3660
